$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '91.133.46', '  +1.89%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.158.63', '  +2.07%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.32%  ')
    ,@(5, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '217.21', '  +1.95%  ')
    ,@(6, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '626.44', '  +0.64%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '1.04', '  +26.04%  ')
    ,@(8, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.373', '  +0.09%  ')
    ,@(9, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.02%  ')
    ,@(10, 'LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '3.151.62', '  +1.98%  ')
    ,@(11, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.749', '  +22.57%  ')
    ,@(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.198', '  +8.91%  ')
    ,@(13, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000250', '  +3.47%  ')
    ,@(14, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '34.80', '  +7.92%  ')
    ,@(15, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '5.55', '  +4.50%  ')
    ,@(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '90.944.54', '  +1.99%  ')
    ,@(17, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.714.75', '  +1.35%  ')
    ,@(18, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.137.53', '  +0.65%  ')
    ,@(19, 'SuiNetwork', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui', '3.74', '  +9.99%  ')
    ,@(20, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0000218', '  +2.15%  ')
    ,@(21, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '14.28', '  +5.96%  ')
    ,@(22, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '446.30', '  +4.84%  ')
    ,@(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '8.93', '  +7.83%  ')
    ,@(24, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.23', '  +5.57%  ')
    ,@(25, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '6.16', '  +11.03%  ')
    ,@(26, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '89.12', '  +6.32%  ')
    ,@(27, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '12.36', '  +2.38%  ')
    ,@(28, 'WrappedeETH', 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth', '3.297.45', '  +1.35%  ')
    ,@(29, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  -0.02%  ')
    ,@(30, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.164', '  -0.11%  ')
    ,@(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '9.20', '  +13.03%  ')
    ,@(32, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '530.18', '  +3.97%  ')
    ,@(33, 'Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '0.906', '  -10.57%  ')
    ,@(34, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '25.35', '  +13.00%  ')
    ,@(35, 'dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '3.73', '  +1.54%  ')
    ,@(36, 'RenderToken', 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render', '7.06', '  +5.28%  ')
    ,@(37, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.143', '  +11.23%  ')
    ,@(38, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.90', '  +5.25%  ')
    ,@(39, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '1.31', '  +4.19%  ')
    ,@(40, 'WhiteBITCoin', 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt', '22.26', '  -0.05%  ')
    ,@(41, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  -0.10%  ')
    ,@(42, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.157', '  +15.76%  ')
    ,@(43, 'PolygonEcosystemToken', 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol', '0.406', '  +11.64%  ')
    ,@(44, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0814', '  +15.93%  ')
    ,@(45, 'USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '1.00', '  +0.02%  ')
    ,@(46, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.94', '  +5.75%  ')
    ,@(47, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '149.21', '  +2.63%  ')
    ,@(48, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.34', '  +9.74%  ')
    ,@(49, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '44.27', '  +1.93%  ')
    ,@(50, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.38', '  +10.22%  ')
    ,@(51, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '172.20', '  +7.67%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Output "Updated $($data.Count) rows"
